# Progress update as of 04-Nov-2025:
#  - Column H ("PERIOD TO EXPIRE") decreases by 1 day for each training row
#  - Column I ("LAST UPDATE") moves from 03-Nov-2025 to 04-Nov-2025
# Applies to rows 3-19 of the "Training Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$firstRow = 3
$lastRow  = 19

$iRange = $ws.Range("I$firstRow`:I$lastRow")

# --- Column H: decrement each numeric value by 1 ---
for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 8)
    $cell.Value = $cell.Value2 - 1
}

# --- Column I: set the new "last update" date text ---
# These cells store the date as literal text (not a real date value), so
# force text formatting first to stop Excel's automatic date recognition
# from turning the string into a date serial number.
$iRange.NumberFormat = "@"
for ($row = $firstRow; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 9).Value = "04-Nov-2025"
}

# Restore the original (General/bordered) cell formatting that the
# NumberFormat change above disturbed, by copying the format from a
# neighboring column that still carries the untouched original style.
$ws.Range("A$firstRow`:A$lastRow").Copy()
$iRange.PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false
